# Initialise learning input automation
# Adds a new "wind_ppl" technology column (E) with learning-curve related
# parameters (learning_par, eos_par, nbr_unit_ref, u_ref) plus a new
# "size" dimension/unit lookup block, and renames the existing
# learning_rate/eos_rate rows to learning_par/eos_par.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: learning_rate -> learning_par, keep dims, add value for wind_ppl
$ws.Range("A22").Value = "learning_par"
$ws.Range("B22").Value = "technology"

# Row 23: eos_rate -> eos_par, keep dims, add value for wind_ppl
$ws.Range("A23").Value = "eos_par"
$ws.Range("B23").Value = "technology"

# New technology column header (E1)
$ws.Range("E1").Value = "wind_ppl"

$ws.Range("E22").Value = 0.9
$ws.Range("E23").Value = 0.9

# Row 24: nbr_unit_ref (number of reference units), unit GW (flagged red)
$ws.Range("A24").Value = "nbr_unit_ref"
$ws.Range("B24").Value = "technology"
$ws.Range("C24").Value = "GW"
$ws.Range("C24").Font.Color = 255
$ws.Range("E24").Value = 100

# Row 25: u_ref (reference unit size), unit GW (flagged red)
$ws.Range("A25").Value = "u_ref"
$ws.Range("B25").Value = "technology"
$ws.Range("C25").Value = "GW"
$ws.Range("C25").Font.Color = 255
$ws.Range("E25").Value = 5

# Row 26: size dimension lookup
$ws.Range("A26").Value = "size"
$ws.Range("E26").Value = "small,medium,large"

# Row 27: u (unit sizes per size), unit GW (flagged red)
$ws.Range("A27").Value = "u"
$ws.Range("C27").Value = "GW"
$ws.Range("C27").Font.Color = 255
$ws.Range("E27").Value = "5,10,50"
$ws.Range("B27").Value = "technology,size"

# Match the page setup recorded by the author's Excel session
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the selection where the author left it
[void]$ws.Range("E27").Select()
